# Apply the edit described by the diff: add a new column C ("sPrice1")
# with a constant value of 1000 in C2, and each subsequent cell (C3:C12)
# referencing the cell directly above it via a simple formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C
$ws.Range("C1").Value = "sPrice1"

# Seed value
$ws.Range("C2").Value = 1000

# Each following row just copies the cell above via formula
$ws.Range("C3").Formula  = "=C2"
$ws.Range("C4").Formula  = "=C3"
$ws.Range("C5").Formula  = "=C4"
$ws.Range("C6").Formula  = "=C5"
$ws.Range("C7").Formula  = "=C6"
$ws.Range("C8").Formula  = "=C7"
$ws.Range("C9").Formula  = "=C8"
$ws.Range("C10").Formula = "=C9"
$ws.Range("C11").Formula = "=C10"
$ws.Range("C12").Formula = "=C11"

# Row heights for the data rows shrink slightly once the third column is in use
for ($r = 3; $r -le 12; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.8
}

# Move the active selection to the newly filled cell
$ws.Range("C2").Select() | Out-Null
